# Applies the price/volume refresh (and the Aave/BabyDogeCoin row swap)
# described by the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '26.390.84'
$ws.Range('E2').Value = '  +0.30%  '
# Row 3
$ws.Range('D3').Value = '1.693.93'
$ws.Range('E3').Value = '  +0.20%  '
# Row 4
$ws.Range('E4').Value = '  +0.27%  '
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.22'
$ws.Range('E5').Value = '  +0.21%  '
# Row 6
$ws.Range('E6').Value = '  +4.62%  '
# Row 7
$ws.Range('E7').Value = '  +0.25%  '
# Row 8
$ws.Range('E8').Value = '  +1.32%  '
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06468'
$ws.Range('E9').Value = '  +0.58%  '
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.02'
$ws.Range('E10').Value = '  -0.10%  '
# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07674'
$ws.Range('E11').Value = '  +2.40%  '
# Row 12
$ws.Range('D12').Value = '1.692.72'
$ws.Range('E12').Value = '  -0.48%  '
# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.543'
$ws.Range('E13').Value = '  -0.49%  '
# Row 14
$ws.Range('E14').Value = '  -0.35%  '
# Row 15
$ws.Range('E15').Value = '  -1.96%  '
# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.49'
$ws.Range('E16').Value = '  +1.41%  '
# Row 17
$ws.Range('D17').Value = '26.425.53'
$ws.Range('E17').Value = '  +0.22%  '
# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.939'
$ws.Range('E18').Value = '  -0.23%  '
# Row 19
$ws.Range('E19').Value = '  +0.31%  '
# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.98'
$ws.Range('E20').Value = '  +0.71%  '
# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '192.09'
$ws.Range('E21').Value = '  +1.18%  '
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.253'
$ws.Range('E22').Value = '  +0.54%  '
# Row 23
$ws.Range('E23').Value = '  +0.26%  '
# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '148.76'
$ws.Range('E24').Value = '  +2.92%  '
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1327'
$ws.Range('E25').Value = '  +7.78%  '
# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.919'
$ws.Range('E26').Value = '  +2.87%  '
# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.77'
$ws.Range('E27').Value = '  -0.67%  '
# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06274'
$ws.Range('E28').Value = '  -5.96%  '
# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.385'
$ws.Range('E29').Value = '  +1.99%  '
# Row 30
$ws.Range('E30').Value = '  +0.10%  '
# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.600'
$ws.Range('E31').Value = '  +0.40%  '
# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.608'
$ws.Range('E32').Value = '  +0.79%  '
# Row 33
$ws.Range('E33').Value = '  +0.79%  '
# Row 34
$ws.Range('E34').Value = '  +1.13%  '
# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6152'
$ws.Range('E35').Value = '  -1.37%  '
# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.414'
$ws.Range('E36').Value = '  +0.75%  '
# Row 37
$ws.Range('E37').Value = '  +0.07%  '
# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01644'
$ws.Range('E38').Value = '  +1.22%  '
# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.181'
$ws.Range('E39').Value = '  -3.22%  '
# Row 40
$ws.Range('D40').Value = '1.117.70'
$ws.Range('E40').Value = '  +0.57%  '
# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8857'
$ws.Range('E41').Value = '  -0.28%  '
# Row 42
$ws.Range('E42').Value = '  -0.14%  '
# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '101.83'
$ws.Range('E43').Value = '  +1.02%  '
# Row 44
$ws.Range('D44').Value = '1.846.03'
$ws.Range('E44').Value = '  +0.40%  '
# Row 45
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '57.49'
$ws.Range('E45').Value = '  +0.85%  '
# Row 46
$ws.Range('B46').Value = 'BabyDogeCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.00000000107'
$ws.Range('E46').Value = '  -4.54%  '
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.012'
$ws.Range('E47').Value = '  +0.09%  '
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.184'
$ws.Range('E48').Value = '  +0.30%  '
# Row 49
$ws.Range('E49').Value = '  +0.30%  '
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.113'
$ws.Range('E50').Value = '  +0.83%  '
# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4304'
$ws.Range('E51').Value = '  -0.02%  '

# Cells that had to be forced to Text format above get their style reset
# back to Normal (General) so no stray number-format styling is left
# behind on cells that were plain/default-styled before the edit.
$ws.Range('D5').Style = 'Normal'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Style = 'Normal'
$ws.Range('D13').Style = 'Normal'
$ws.Range('D16').Style = 'Normal'
$ws.Range('D18').Style = 'Normal'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D24').Style = 'Normal'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').Style = 'Normal'
$ws.Range('D29').Style = 'Normal'
$ws.Range('D31').Style = 'Normal'
$ws.Range('D32').Style = 'Normal'
$ws.Range('D35').Style = 'Normal'
$ws.Range('D36').Style = 'Normal'
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Style = 'Normal'
$ws.Range('D41').Style = 'Normal'
$ws.Range('D43').Style = 'Normal'
$ws.Range('D45').Style = 'Normal'
$ws.Range('D46').Style = 'Normal'
$ws.Range('D47').Style = 'Normal'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D50').Style = 'Normal'
$ws.Range('D51').Style = 'Normal'
